$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (these cells currently hold numbers
# in the "before" workbook; they must become text values instead).
$values = [ordered]@{
    "H2"  = "2017"
    "H3"  = "2013"
    "H4"  = "2013"
    "H5"  = "2018-04-27 00:00:00"
    "H6"  = "2016-07-21 00:00:00"
    "H7"  = "2012-09-27 00:00:00"
    "H8"  = "2020"
    "H9"  = "2019"
    "H10" = "2020"
    "H11" = "201409"
    "H12" = "201808"
    "H13" = "201910"
    "H14" = "2001-10-01 00:00:00"
    "H15" = "2001-01-01 00:00:00"
    "H16" = "2015-09-10 00:00:00"
}

foreach ($addr in $values.Keys) {
    $cell = $ws.Range($addr)
    # Force a text number format before assigning, so the numeric-looking
    # string is stored as text rather than being reinterpreted as a number.
    $cell.NumberFormat = "@"
    $cell.Value2 = $values[$addr]
    # Re-apply the default "Normal" style so no left-over per-cell style
    # (and in particular no date number format) remains on the cell.
    $cell.Style = "Normal"
}

$wb.Save()
